$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data (A14:D15)
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = 45441
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 3

$ws.Range("A15").Value = 12
$ws.Range("B15").Value = 45448
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 3

# Update the selection to A16 like in the diff
$ws.Range("A16").Select()

# Update the chart series to extend their source ranges to row 15
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Planilha1!`$C`$3,Planilha1!`$B`$4:`$B`$15,Planilha1!`$C`$4:`$C`$15,1)"

$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Planilha1!`$D`$3,Planilha1!`$B`$4:`$B`$15,Planilha1!`$D`$4:`$D`$15,2)"
